# Fix the example to refer to strains by their correct name:
#   JJS-MGP1  -> JJS-MGP001
#   JJS-MGP20 -> JJS-MGP020
# These strain names appear in the "openbis-data" sheet, column A (rows 2-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-data")

for ($r = 2; $r -le 8; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -eq "JJS-MGP1") {
        $cell.Value = "JJS-MGP001"
    } elseif ($val -eq "JJS-MGP20") {
        $cell.Value = "JJS-MGP020"
    }
}

# Restore the sheet's previous selection (it moved to B26 in the saved file).
$ws.Activate()
$ws.Range("B26").Select()
